$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency "Price" values in column D.
# The leading apostrophe forces Excel to store the value as text
# (preserving exact numeric formatting, e.g. trailing zeros), matching
# the original inline-string cell contents. Resetting the style back to
# "Normal" avoids Excel applying an implicit @ (Text) number format to
# the cell, keeping the cell style identical to the source workbook.
$ws.Range("D2").Value = "'242.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.415"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05888"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.437"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.539"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8103"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9305"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07424"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03369"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03043"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09349"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.844"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04669"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005923"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.005891"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001265"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004900"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00006799"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.564"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").Value = "'0.0002297"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03975"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006184"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002569"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.009717"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005182"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6704"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002389"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
